$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 6 ("Estender implementação da MM e da soma concorrente para qualquer n") as done
$ws.Range("C6").Value = "FEITO"

# Update the active selection to C7
$ws.Range("C7").Select()
